$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (BNB)
$ws.Range("D2").Value = "'245.82"

# Row 3 (OKB)
$ws.Range("D3").Value = "'23.80"

# Row 4 (HuobiToken)
$ws.Range("D4").Value = "'5.351"

# Row 5 (Cronos)
$ws.Range("D5").Value = "'0.05840"

# Row 6 (KuCoinToken)
$ws.Range("D6").Value = "'6.483"

# Row 7 (GateToken)
$ws.Range("D7").Value = "'3.359"

# Row 8 (MXToken)
$ws.Range("D8").Value = "'0.8121"

# Row 9 (FTXToken)
$ws.Range("D9").Value = "'0.9243"

# Row 10: was One -> now WazirX
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1413"
$ws.Range("E10").Value = "9WazirXWRX"

# Row 11: was WazirX -> now MandalaExchangeToken
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07366"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

# Row 12: was MandalaExchangeToken -> now LiechtensteinCryptoassetsExchange
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03107"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

# Row 13: was LiechtensteinCryptoassetsExchange -> now BitrueCoin
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03063"
$ws.Range("E13").Value = "12BitrueCoinBTR"

# Row 14: was BitrueCoin -> now BitMartToken
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09389"
$ws.Range("E14").Value = "13BitMartTokenBMX"

# Row 15: was BitMartToken -> now MCDex
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.859"
$ws.Range("E15").Value = "14MCDexMCB"

# Row 16: was MCDex -> now BitForexToken
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001557"
$ws.Range("E16").Value = "15BitForexTokenBF"

# Row 17: was BitForexToken -> now CoinExToken
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04710"
$ws.Range("E17").Value = "16CoinExTokenCET"

# Row 18: was CoinExToken -> now One
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0006068"
$ws.Range("E18").Value = "17OneONE"

# Row 20 (BitKan)
$ws.Range("D20").Value = "'0.001250"

# Row 21 (HotbitToken)
$ws.Range("D21").Value = "'0.004689"

# Row 22 (NitroEx)
$ws.Range("D22").Value = "'0.00008820"
$ws.Range("E22").Value = "21NitroExNTXBestin24h"

# Row 23 (LEO)
$ws.Range("D23").Value = "'3.595"

# Row 25 (BitpandaEcosystemToken)
$ws.Range("D25").Value = "'0.3228"

# Row 28 (UpBots)
$ws.Range("D28").Value = "'0.0002661"

# Row 40 (IDEX)
$ws.Range("D40").Value = "'0.03846"

# Row 41 (KickToken)
$ws.Range("D41").Value = "'0.006424"

# Row 42 (BKEXToken)
$ws.Range("D42").Value = "'0.1067"

# Row 43 (CEJI)
$ws.Range("D43").Value = "'0.002947"

# Row 44 (LocalTraders)
$ws.Range("D44").Value = "'0.008576"

# Row 45 (CoinLion)
$ws.Range("D45").Value = "'0.00005269"

# Row 47 (CoinbaseStockToken)
$ws.Range("D47").Value = "'0.6540"

# Row 48 (BOLO)
$ws.Range("D48").Value = "'0.001863"

# Row 49 (CryptobidCoin)
$ws.Range("D49").Value = "'0.00002106"

# Row 50 (SpecialPowerGold)
$ws.Range("D50").Value = "'0.0002006"
